$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("B2").Value = 0.5280971303754995
$ws.Range("C2").Value = 0.2284544404075213
$ws.Range("D2").Value = 0.06097336337137449
$ws.Range("E2").Value = 0.1388871599330486
$ws.Range("F2").Value = 1.181426164879099
$ws.Range("K2").Value = 0.2619865053838168
$ws.Range("L2").Value = 0.1918557080101237
$ws.Range("M2").Value = 0.1544339646903232
$ws.Range("O2").Value = 4.300780462613432
$ws.Range("B3").Value = 0.4935342816666264
$ws.Range("C3").Value = 0.2284731720951072
$ws.Range("D3").Value = 0.05921214971377253
$ws.Range("E3").Value = 0.1394203463867356
$ws.Range("F3").Value = 1.18256094622889
$ws.Range("K3").Value = 0.2308765021766419
$ws.Range("L3").Value = 0.1893104981596068
$ws.Range("M3").Value = 0.1479391547985891
$ws.Range("O3").Value = 4.319256954297515
$ws.Range("B4").Value = 0.4724925182315189
$ws.Range("C4").Value = 0.2285009141132477
$ws.Range("D4").Value = 0.05811692019202752
$ws.Range("E4").Value = 0.139799732201686
$ws.Range("F4").Value = 1.183865847586276
$ws.Range("K4").Value = 0.2117681803738805
$ws.Range("L4").Value = 0.1878389323735092
$ws.Range("M4").Value = 0.1440178982112599
$ws.Range("O4").Value = 4.332642878461513
$ws.Range("B5").Value = 0.463963627350239
$ws.Range("C5").Value = 0.2285163272789212
$ws.Range("D5").Value = 0.05766714504215287
$ws.Range("E5").Value = 0.1399674318446351
$ws.Range("F5").Value = 1.184550604753554
$ws.Range("K5").Value = 0.2039801166259991
$ws.Range("L5").Value = 0.1872622499611367
$ws.Range("M5").Value = 0.1424368054558194
$ws.Range("O5").Value = 4.338611232546015
$ws.Range("B6").Value = 0.4625501937119907
$ws.Range("C6").Value = 0.2285191354575176
$ws.Range("D6").Value = 0.05759225180126037
$ws.Range("E6").Value = 0.1399960698248606
$ws.Range("F6").Value = 1.184673551213635
$ws.Range("K6").Value = 0.2026868500055201
$ws.Range("L6").Value = 0.1871678827135739
$ws.Range("M6").Value = 0.1421752871704349
$ws.Range("O6").Value = 4.339633291634613
$ws.Range("B7").Value = 0.4723773084657807
$ws.Range("C7").Value = 0.2285011053134376
$ws.Range("D7").Value = 0.05811086835249313
$ws.Range("E7").Value = 0.139801940803407
$ws.Range("F7").Value = 1.183874462871707
$ws.Range("K7").Value = 0.2116631524149568
$ws.Range("L7").Value = 0.187831061858283
$ws.Range("M7").Value = 0.1439965066677296
$ws.Range("O7").Value = 4.332721290467646
$ws.Range("B8").Value = 0.5161428231544392
$ws.Range("C8").Value = 0.2284575427009905
$ws.Range("D8").Value = 0.06036897987339529
$ws.Range("E8").Value = 0.1390602199417277
$ws.Range("F8").Value = 1.18169126153736
$ws.Range("K8").Value = 0.251261390787505
$ws.Range("L8").Value = 0.1909592243999114
$ws.Range("M8").Value = 0.1521807975338625
$ws.Range("O8").Value = 4.306727593939883
$ws.Range("B9").Value = 0.6033755538961429
$ws.Range("C9").Value = 0.2284999591784072
$ws.Range("D9").Value = 0.064686676617562
$ws.Range("E9").Value = 0.1380175735434328
$ws.Range("F9").Value = 1.182233243122475
$ws.Range("K9").Value = 0.3288466244278538
$ws.Range("L9").Value = 0.1978151420176815
$ws.Range("M9").Value = 0.168754667497506
$ws.Range("O9").Value = 4.271945276234874
$ws.Range("B10").Value = 0.6683052358698092
$ws.Range("C10").Value = 0.2286077519747849
$ws.Range("D10").Value = 0.06779093721643648
$ws.Range("E10").Value = 0.1375016505422195
$ws.Range("F10").Value = 1.185570924565738
$ws.Range("K10").Value = 0.3857946177257361
$ws.Range("L10").Value = 0.2032901012923531
$ws.Range("M10").Value = 0.1812475312341064
$ws.Range("O10").Value = 4.256257972840871
$ws.Range("B11").Value = 0.6980221968193518
$ws.Range("C11").Value = 0.2286731601055649
$ws.Range("D11").Value = 0.06918827726053678
$ws.Range("E11").Value = 0.1373210524912238
$ws.Range("F11").Value = 1.18772754264468
$ws.Range("K11").Value = 0.4116876050145777
$ws.Range("L11").Value = 0.2058755181214167
$ws.Range("M11").Value = 0.1869987861070328
$ws.Range("O11").Value = 4.251263587692051
$ws.Range("B12").Value = 0.7093006842715397
$ws.Range("C12").Value = 0.2287002600376695
$ws.Range("D12").Value = 0.06971526896246161
$ws.Range("E12").Value = 0.1372604262444881
$ws.Range("F12").Value = 1.18863594090773
$ws.Range("K12").Value = 0.4214904134902895
$ws.Range("L12").Value = 0.2068681368366185
$ws.Range("M12").Value = 0.1891863525037394
$ws.Range("O12").Value = 4.249680246584717
$ws.Range("B13").Value = 0.7068705433748335
$ws.Range("C13").Value = 0.2286943202615603
$ws.Range("D13").Value = 0.06960186786767508
$ws.Range("E13").Value = 0.1372731382294354
$ws.Range("F13").Value = 1.188436222033999
$ws.Range("K13").Value = 0.4193793127943479
$ws.Range("L13").Value = 0.2066537557147541
$ws.Range("M13").Value = 0.1887147917506482
$ws.Range("O13").Value = 4.250007553553871
$ws.Range("B14").Value = 0.6989495806639923
$ws.Range("C14").Value = 0.2286753430418855
$ws.Range("D14").Value = 0.06923167641232197
$ws.Range("E14").Value = 0.137315909261833
$ws.Range("F14").Value = 1.187800438933238
$ws.Range("K14").Value = 0.41249413667299
$ws.Range("L14").Value = 0.2059569096861935
$ws.Range("M14").Value = 0.1871785649440554
$ws.Range("O14").Value = 4.251127154418242
$ws.Range("B15").Value = 0.6941010423527985
$ws.Range("C15").Value = 0.2286640218695055
$ws.Range("D15").Value = 0.06900464273044093
$ws.Range("E15").Value = 0.1373431181354938
$ws.Range("F15").Value = 1.187422948328418
$ws.Range("K15").Value = 0.4082764543381359
$ws.Range("L15").Value = 0.2055318378076549
$ws.Range("M15").Value = 0.1862388403272845
$ws.Range("O15").Value = 4.251853040963908
$ws.Range("B16").Value = 0.6663667373215389
$ws.Range("C16").Value = 0.2286038048804073
$ws.Range("D16").Value = 0.06769931799697559
$ws.Range("E16").Value = 0.1375145401930737
$ws.Range("F16").Value = 1.185442824147657
$ws.Range("K16").Value = 0.3841021518125842
$ws.Range("L16").Value = 0.2031230414081762
$ws.Range("M16").Value = 0.1808730355445221
$ws.Range("O16").Value = 4.256627461071446
$ws.Range("B17").Value = 0.6493983586539684
$ws.Range("C17").Value = 0.2285710406287436
$ws.Range("D17").Value = 0.06689473637174359
$ws.Range("E17").Value = 0.1376335458175753
$ws.Range("F17").Value = 1.184391527144925
$ws.Range("K17").Value = 0.3692683701772808
$ws.Range("L17").Value = 0.2016695691992823
$ws.Range("M17").Value = 0.1775986713558879
$ws.Range("O17").Value = 4.260104941671898
$ws.Range("B18").Value = 0.639655584783128
$ws.Range("C18").Value = 0.2285537379304543
$ws.Range("D18").Value = 0.06643057021217658
$ws.Range("E18").Value = 0.1377070876662625
$ws.Range("F18").Value = 1.183846930933896
$ws.Range("K18").Value = 0.3607351796931653
$ws.Range("L18").Value = 0.2008424998466012
$ws.Range("M18").Value = 0.1757217687344621
$ws.Range("O18").Value = 4.262306701687351
$ws.Range("B19").Value = 0.6363597845483469
$ws.Range("C19").Value = 0.2285481451563953
$ws.Range("D19").Value = 0.06627317310176295
$ws.Range("E19").Value = 0.1377328629085603
$ws.Range("F19").Value = 1.183672860715461
$ws.Range("K19").Value = 0.3578457968743578
$ws.Range("L19").Value = 0.2005640038084238
$ws.Range("M19").Value = 0.1750873887548963
$ws.Range("O19").Value = 4.263086807829353
$ws.Range("B20").Value = 0.6512029174821521
$ws.Range("C20").Value = 0.2285743689717918
$ws.Range("D20").Value = 0.06698052969363744
$ws.Range("E20").Value = 0.1376203504858644
$ws.Range("F20").Value = 1.18449722161094
$ws.Range("K20").Value = 0.3708475780983065
$ws.Range("L20").Value = 0.2018233700639911
$ws.Range("M20").Value = 0.1779465688739634
$ws.Range("O20").Value = 4.259713893052805
$ws.Range("B21").Value = 0.7012754768942386
$ws.Range("C21").Value = 0.2286808540230396
$ws.Range("D21").Value = 0.06934046915841208
$ws.Range("E21").Value = 0.1373031358404013
$ws.Range("F21").Value = 1.187984694628156
$ws.Range("K21").Value = 0.4145165458492386
$ws.Range("L21").Value = 0.2061612222383786
$ws.Range("M21").Value = 0.1876295298963058
$ws.Range("O21").Value = 4.250789944102905
$ws.Range("B22").Value = 0.7341480999825194
$ws.Range("C22").Value = 0.2287640262117634
$ws.Range("D22").Value = 0.07087028030134945
$ws.Range("E22").Value = 0.1371410546165244
$ws.Range("F22").Value = 1.190798630028738
$ws.Range("K22").Value = 0.4430431074109435
$ws.Range("L22").Value = 0.2090753718649836
$ws.Range("M22").Value = 0.1940143382060455
$ws.Range("O22").Value = 4.246752411474347
$ws.Range("B23").Value = 0.7165900754359313
$ws.Range("C23").Value = 0.2287184008081908
$ws.Range("D23").Value = 0.07005494645217425
$ws.Range("E23").Value = 0.1372234268173891
$ws.Range("F23").Value = 1.189247874675345
$ws.Range("K23").Value = 0.4278193381062749
$ws.Range("L23").Value = 0.2075128160856252
$ws.Range("M23").Value = 0.1906015199495386
$ws.Range("O23").Value = 4.248743119533685
$ws.Range("B24").Value = 0.6503870372294784
$ws.Range("C24").Value = 0.2285728594479082
$ws.Range("D24").Value = 0.06694174751910964
$ws.Range("E24").Value = 0.1376263001336877
$ws.Range("F24").Value = 1.184449250838483
$ws.Range("K24").Value = 0.3701336339010197
$ws.Range("L24").Value = 0.2017538100531482
$ws.Range("M24").Value = 0.177789267038726
$ws.Range("O24").Value = 4.259890055314429
$ws.Range("B25").Value = 0.5796279470778245
$ws.Range("C25").Value = 0.2284749341374592
$ws.Range("D25").Value = 0.06353051403426235
$ws.Range("E25").Value = 0.1382556448363506
$ws.Range("F25").Value = 1.181570350968279
$ws.Range("K25").Value = 0.3078662799447613
$ws.Range("L25").Value = 0.1958833880421054
$ws.Range("M25").Value = 0.1642152292169747
$ws.Range("O25").Value = 4.279621775039288
